# Scheduled Sheets refresh: pushes newly-recalculated market/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns, i.e. H:N) into the
# per-job worksheets. Generated by the runner from the latest price pull.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of {cell, value} edits to apply.
$sheetEdits = @{
    "ALC" = @(
        @{Cell="H28"; Val=44443.13},
        @{Cell="I28"; Val=56713.332},
        @{Cell="J28"; Val=270.4},
        @{Cell="K28"; Val=56713.332},
        @{Cell="L28"; Val=270.4},
        @{Cell="M28"; Val=-56228.332},
        @{Cell="N28"; Val=-1240.4},
        @{Cell="H55"; Val=125859},
        @{Cell="I55"; Val=32.4},
        @{Cell="J55"; Val=335570},
        @{Cell="K55"; Val=32.4},
        @{Cell="L55"; Val=335570},
        @{Cell="M55"; Val=181.6},
        @{Cell="N55"; Val=-335998},
        @{Cell="H106"; Val=2862.25},
        @{Cell="I106"; Val=3349.6667},
        @{Cell="K106"; Val=3349.6667},
        @{Cell="M106"; Val=-2718.6667},
        @{Cell="H111"; Val=62966.234},
        @{Cell="I111"; Val=103924.4},
        @{Cell="J111"; Val=4454.5713},
        @{Cell="K111"; Val=311773.2},
        @{Cell="L111"; Val=13363.7139},
        @{Cell="M111"; Val=-308706.2},
        @{Cell="N111"; Val=-19497.7139},
        @{Cell="H118"; Val=746},
        @{Cell="I118"; Val=592.1111},
        @{Cell="J118"; Val=1207.6666},
        @{Cell="K118"; Val=1776.3333},
        @{Cell="L118"; Val=3622.9998},
        @{Cell="M118"; Val=-119.3332999999998},
        @{Cell="N118"; Val=-6936.9998},
        @{Cell="H125"; Val=6176177.5},
        @{Cell="I125"; Val=2121.8572},
        @{Cell="J125"; Val=10105122},
        @{Cell="K125"; Val=19096.7148},
        @{Cell="L125"; Val=90946098},
        @{Cell="M125"; Val=-16636.7148},
        @{Cell="N125"; Val=-90951018},
        @{Cell="H132"; Val=3879.92},
        @{Cell="I132"; Val=4397.1333},
        @{Cell="K132"; Val=13191.3999},
        @{Cell="M132"; Val=-10661.3999},
        @{Cell="H138"; Val=5885.2646},
        @{Cell="J138"; Val=7352.1665},
        @{Cell="L138"; Val=22056.4995},
        @{Cell="N138"; Val=-32336.4995},
    )
    "ARM" = @(
        @{Cell="H32"; Val=4065.9016},
        @{Cell="I32"; Val=4065.9016},
        @{Cell="K32"; Val=4065.9016},
        @{Cell="M32"; Val=-3778.9016},
        @{Cell="H45"; Val=2432.158},
        @{Cell="I45"; Val=1715.2142},
        @{Cell="J45"; Val=4439.6},
        @{Cell="K45"; Val=1715.2142},
        @{Cell="L45"; Val=4439.6},
        @{Cell="M45"; Val=-1338.2142},
        @{Cell="N45"; Val=-5193.6},
        @{Cell="H74"; Val=1967.5454},
        @{Cell="I74"; Val=1646.6842},
        @{Cell="J74"; Val=3999.6667},
        @{Cell="K74"; Val=1646.6842},
        @{Cell="L74"; Val=3999.6667},
        @{Cell="M74"; Val=-772.6841999999999},
        @{Cell="N74"; Val=-5747.6667},
        @{Cell="H77"; Val=1967.5454},
        @{Cell="I77"; Val=1646.6842},
        @{Cell="J77"; Val=3999.6667},
        @{Cell="K77"; Val=8233.421},
        @{Cell="L77"; Val=19998.3335},
        @{Cell="M77"; Val=-3865.421},
        @{Cell="N77"; Val=-28734.3335},
        @{Cell="H110"; Val=148457.23},
        @{Cell="I110"; Val=148457.23},
        @{Cell="K110"; Val=148457.23},
        @{Cell="M110"; Val=-146412.23},
        @{Cell="H122"; Val=3290.8},
        @{Cell="I122"; Val=1855.28},
        @{Cell="J122"; Val=5683.3335},
        @{Cell="K122"; Val=5565.84},
        @{Cell="L122"; Val=17050.0005},
        @{Cell="M122"; Val=-3115.84},
        @{Cell="N122"; Val=-21950.0005},
    )
    "BSM" = @(
        @{Cell="H134"; Val=30357.805},
        @{Cell="I134"; Val=1829.3667},
        @{Cell="K134"; Val=5488.1001},
        @{Cell="M134"; Val=-2953.1001},
        @{Cell="H139"; Val=81000},
        @{Cell="J139"; Val=81000},
        @{Cell="L139"; Val=81000},
        @{Cell="N139"; Val=-91280},
    )
    "CRP" = @(
        @{Cell="H16"; Val=1390.6666},
        @{Cell="I16"; Val=1308.2},
        @{Cell="J16"; Val=1493.75},
        @{Cell="K16"; Val=1308.2},
        @{Cell="L16"; Val=1493.75},
        @{Cell="M16"; Val=-1021.2},
        @{Cell="N16"; Val=-2067.75},
        @{Cell="H113"; Val=1390.6666},
        @{Cell="I113"; Val=1308.2},
        @{Cell="J113"; Val=1493.75},
        @{Cell="K113"; Val=1308.2},
        @{Cell="L113"; Val=1493.75},
        @{Cell="M113"; Val=861.8},
        @{Cell="N113"; Val=-5833.75},
        @{Cell="H134"; Val=558342},
        @{Cell="I134"; Val=359620.7},
        @{Cell="K134"; Val=1078862.1},
        @{Cell="M134"; Val=-1076327.1},
    )
    "CUL" = @(
        @{Cell="H113"; Val=1324119.9},
        @{Cell="J113"; Val=1409.3334},
        @{Cell="L113"; Val=4228.0002},
        @{Cell="N113"; Val=-8568.0002},
        @{Cell="H116"; Val=1674},
        @{Cell="J116"; Val=1674},
        @{Cell="L116"; Val=5022},
        @{Cell="N116"; Val=-11906},
        @{Cell="H137"; Val=2579.5},
        @{Cell="I137"; Val=1428.6666},
        @{Cell="K137"; Val=4285.9998},
        @{Cell="M137"; Val=814.0002000000004},
        @{Cell="H139"; Val=5096.4585},
        @{Cell="I139"; Val=4333.8945},
        @{Cell="K139"; Val=13001.6835},
        @{Cell="M139"; Val=-7861.683500000001},
        @{Cell="H140"; Val=2340},
        @{Cell="I140"; Val=2340},
        @{Cell="K140"; Val=7020},
        @{Cell="M140"; Val=-1840},
    )
    "GSM" = @(
        @{Cell="H36"; Val=4600},
        @{Cell="J36"; Val=4600},
        @{Cell="L36"; Val=4600},
        @{Cell="N36"; Val=-5570},
        @{Cell="H70"; Val=100008130},
        @{Cell="I70"; Val=10056},
        @{Cell="J70"; Val=333336960},
        @{Cell="K70"; Val=10056},
        @{Cell="L70"; Val=333336960},
        @{Cell="M70"; Val=-9786},
        @{Cell="N70"; Val=-333337500},
        @{Cell="H73"; Val=100008130},
        @{Cell="I73"; Val=10056},
        @{Cell="J73"; Val=333336960},
        @{Cell="K73"; Val=10056},
        @{Cell="L73"; Val=333336960},
        @{Cell="M73"; Val=-9120},
        @{Cell="N73"; Val=-333338832},
        @{Cell="H113"; Val=407357.6},
        @{Cell="I113"; Val=717077.0600000001},
        @{Cell="J113"; Val=13169.182},
        @{Cell="K113"; Val=717077.0600000001},
        @{Cell="L113"; Val=13169.182},
        @{Cell="M113"; Val=-714907.0600000001},
        @{Cell="N113"; Val=-17509.182},
        @{Cell="H122"; Val=3545.9546},
        @{Cell="I122"; Val=1750},
        @{Cell="J122"; Val=4384.067},
        @{Cell="K122"; Val=5250},
        @{Cell="L122"; Val=13152.201},
        @{Cell="M122"; Val=-2800},
        @{Cell="N122"; Val=-18052.201},
        @{Cell="H132"; Val=617729.75},
        @{Cell="I132"; Val=2011244.6},
        @{Cell="K132"; Val=6033733.800000001},
        @{Cell="M132"; Val=-6031203.800000001},
    )
    "LTW" = @(
        @{Cell="H7"; Val=230765.6},
        @{Cell="I7"; Val=3384.4814},
        @{Cell="J7"; Val=591900.3},
        @{Cell="K7"; Val=3384.4814},
        @{Cell="L7"; Val=591900.3},
        @{Cell="M7"; Val=-3272.4814},
        @{Cell="N7"; Val=-592124.3},
        @{Cell="H46"; Val=4027.4482},
        @{Cell="I46"; Val=5186.6},
        @{Cell="J46"; Val=2785.5},
        @{Cell="K46"; Val=5186.6},
        @{Cell="L46"; Val=2785.5},
        @{Cell="M46"; Val=-4998.6},
        @{Cell="N46"; Val=-3161.5},
        @{Cell="H47"; Val=182208.33},
        @{Cell="I47"; Val=1000000},
        @{Cell="J47"; Val=18650},
        @{Cell="K47"; Val=1000000},
        @{Cell="L47"; Val=18650},
        @{Cell="M47"; Val=-999510},
        @{Cell="N47"; Val=-19630},
        @{Cell="H52"; Val=182208.33},
        @{Cell="I52"; Val=1000000},
        @{Cell="J52"; Val=18650},
        @{Cell="K52"; Val=1000000},
        @{Cell="L52"; Val=18650},
        @{Cell="M52"; Val=-999767},
        @{Cell="N52"; Val=-19116},
        @{Cell="H61"; Val=4318.625},
        @{Cell="I61"; Val=3938.8},
        @{Cell="J61"; Val=4951.6665},
        @{Cell="K61"; Val=3938.8},
        @{Cell="L61"; Val=4951.6665},
        @{Cell="M61"; Val=-3736.8},
        @{Cell="N61"; Val=-5355.6665},
        @{Cell="H70"; Val=21125},
        @{Cell="J70"; Val=21125},
        @{Cell="L70"; Val=21125},
        @{Cell="N70"; Val=-21665},
        @{Cell="H73"; Val=21125},
        @{Cell="J73"; Val=21125},
        @{Cell="L73"; Val=21125},
        @{Cell="N73"; Val=-22997},
        @{Cell="H93"; Val=2007.909},
        @{Cell="I93"; Val=1979.875},
        @{Cell="K93"; Val=1979.875},
        @{Cell="M93"; Val=-731.875},
        @{Cell="H113"; Val=4318.625},
        @{Cell="I113"; Val=3938.8},
        @{Cell="J113"; Val=4951.6665},
        @{Cell="K113"; Val=3938.8},
        @{Cell="L113"; Val=4951.6665},
        @{Cell="M113"; Val=-1768.8},
        @{Cell="N113"; Val=-9291.666499999999},
        @{Cell="H124"; Val=82000},
        @{Cell="J124"; Val=82000},
        @{Cell="L124"; Val=82000},
        @{Cell="N124"; Val=-91820},
        @{Cell="H126"; Val=230765.6},
        @{Cell="I126"; Val=3384.4814},
        @{Cell="J126"; Val=591900.3},
        @{Cell="K126"; Val=10153.4442},
        @{Cell="L126"; Val=1775700.9},
        @{Cell="M126"; Val=-7683.4442},
        @{Cell="N126"; Val=-1780640.9},
    )
    "WVR" = @(
        @{Cell="H81"; Val=1893},
        @{Cell="I81"; Val=1274.4286},
        @{Cell="J81"; Val=2975.5},
        @{Cell="K81"; Val=2548.8572},
        @{Cell="L81"; Val=5951},
        @{Cell="M81"; Val=-1487.8572},
        @{Cell="N81"; Val=-8073},
        @{Cell="H84"; Val=1893},
        @{Cell="I84"; Val=1274.4286},
        @{Cell="J84"; Val=2975.5},
        @{Cell="K84"; Val=12744.286},
        @{Cell="L84"; Val=29755},
        @{Cell="M84"; Val=-7440.286},
        @{Cell="N84"; Val=-40363},
        @{Cell="H113"; Val=261.625},
        @{Cell="I113"; Val=302.375},
        @{Cell="J113"; Val=220.875},
        @{Cell="K113"; Val=907.125},
        @{Cell="L113"; Val=662.625},
        @{Cell="M113"; Val=1262.875},
        @{Cell="N113"; Val=-5002.625},
        @{Cell="H122"; Val=29414094},
        @{Cell="I122"; Val=32259894},
        @{Cell="K122"; Val=96779682},
        @{Cell="M122"; Val=-96777232},
        @{Cell="H126"; Val=4591.5557},
        @{Cell="I126"; Val=3450},
        @{Cell="K126"; Val=10350},
        @{Cell="M126"; Val=-7880},
        @{Cell="H135"; Val=64449.7},
        @{Cell="J135"; Val=64449.7},
        @{Cell="L135"; Val=64449.7},
        @{Cell="N135"; Val=-74589.7},
    )
}

foreach ($sheetName in $sheetEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($edit in $sheetEdits[$sheetName]) {
        $ws.Range($edit.Cell).Value = $edit.Val
    }
}
